$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in Status ("E") column for the rows that were missing it
$ws.Range("E11").Value = "Done - can't repropuce error"
$ws.Range("E12").Value = "Done"
$ws.Range("E14").Value = "Done"

# Update the selected cell shown in the sheet view
$ws.Range("E12").Select()
